# Update cryptos list with latest prices/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value happens to look like a plain number need to stay
# text-typed (matching the source data, which always stores these as text -
# note the "." thousands separators and "%" suffixes elsewhere in the column)
# so pre-format them as Text before writing the values.
$ws.Range("D5,D6,D7,D9,D10,D11,D12,D14,D16,D20,D21,D22,D23,D24,D25,D28,D30,D31,D32,D34,D35,D37,D41,D42,D43,D44,D47,D48,D49").NumberFormat = "@"

$ws.Range("D2").Value = '43.649.01'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '2.283.78'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").Value = '109.73'
$ws.Range("E5").Value = '  +14.26%  '
$ws.Range("D6").Value = '266.85'
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  +1.11%  '
$ws.Range("E8").Value = '  +0.34%  '
$ws.Range("D9").Value = '0.612'
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("D10").Value = '47.19'
$ws.Range("E10").Value = '  +3.48%  '
$ws.Range("D11").Value = '0.0941'
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").Value = '8.78'
$ws.Range("E12").Value = '  +9.92%  '
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("D14").Value = '15.59'
$ws.Range("E14").Value = '  +1.30%  '
$ws.Range("D15").Value = '2.626.44'
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").Value = '0.839'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").Value = '2.279.06'
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").Value = '43.493.91'
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("D20").Value = '6.62'
$ws.Range("E20").Value = '  +6.87%  '
$ws.Range("D21").Value = '72.08'
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").Value = '2.45'
$ws.Range("E22").Value = '  -4.62%  '
$ws.Range("D23").Value = '231.23'
$ws.Range("E23").Value = '  -0.73%  '
$ws.Range("D24").Value = '9.59'
$ws.Range("E24").Value = '  +5.21%  '
$ws.Range("D25").Value = '2.76'
$ws.Range("E25").Value = '  +8.18%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  +3.03%  '
$ws.Range("D28").Value = '41.77'
$ws.Range("E28").Value = '  +4.14%  '
$ws.Range("E29").Value = '  -2.27%  '
$ws.Range("D30").Value = '2.26'
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("D31").Value = '175.79'
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("D32").Value = '21.42'
$ws.Range("E32").Value = '  -2.48%  '
$ws.Range("E33").Value = '  +2.53%  '
$ws.Range("D34").Value = '5.58'
$ws.Range("E34").Value = '  +4.14%  '
$ws.Range("D35").Value = '0.126'
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("E36").Value = '  +7.25%  '
$ws.Range("D37").Value = '0.0359'
$ws.Range("E37").Value = '  +1.83%  '
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("E39").Value = '  +11.73%  '
$ws.Range("E40").Value = '  +3.53%  '
$ws.Range("D41").Value = '0.241'
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").Value = '13.48'
$ws.Range("E42").Value = '  +9.77%  '
$ws.Range("D43").Value = '71.32'
$ws.Range("E43").Value = '  +8.73%  '
$ws.Range("D44").Value = '6.25'
$ws.Range("E44").Value = '  +20.65%  '
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("E46").Value = '  +2.18%  '
$ws.Range("D47").Value = '8.78'
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").Value = '0.0999'
$ws.Range("E48").Value = '  -1.67%  '
$ws.Range("D49").Value = '101.22'
$ws.Range("E49").Value = '  +4.26%  '
$ws.Range("E50").Value = '  +1.69%  '
$ws.Range("E51").Value = '  +4.86%  '

Write-Output "Updated $($ws.Name) with latest crypto prices"
